# Tutorial 6 solution update:
#  - Reformat the Date column (A3:A21) from DD/MM/YYYY to DD-MM-YYYY
#  - Correct the attendance tallies for the first two data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> corrected date text (dashes instead of slashes)
$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force text storage so Excel doesn't reinterpret the dashed string as a date
    # serial (relevant for day<=12 values like "01-08-2022"), then restore the
    # default "Normal" style so no stray number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.Style = "Normal"
}

# Attendance count corrections
$ws.Cells.Item(3, 4).Value = 1   # D3: Total Attendance Count 0 -> 1
$ws.Cells.Item(3, 7).Value = 1   # G3: Invalid 0 -> 1

$ws.Cells.Item(4, 4).Value = 1   # D4: Total Attendance Count 0 -> 1
$ws.Cells.Item(4, 5).Value = 1   # E4: Real 0 -> 1
$ws.Cells.Item(4, 8).Value = 0   # H4: Absent 1 -> 0
